# Weekly update: insert a new price record for "Zapallo italiano"
# (Comercializadora del Agro de Limarí) as row 40, pushing the existing
# rows 40-57 down to 41-58 (dimension grows from A1:R57 to A1:R58).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 40..57 down by one, leaving a blank row 40 to fill in.
$ws.Rows.Item(40).Insert()

# New row 40 values.
$ws.Cells.Item(40, 1).Value  = 2
$ws.Cells.Item(40, 2).Value  = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(40, 3).Value  = "Coquimbo"
$ws.Cells.Item(40, 4).Value  = 44609
$ws.Cells.Item(40, 5).Value  = 4
$ws.Cells.Item(40, 6).Value  = 100112032
$ws.Cells.Item(40, 7).Value  = "Zapallo italiano"
$ws.Cells.Item(40, 8).Value  = "Sin especificar"
$ws.Cells.Item(40, 9).Value  = "Primera"
$ws.Cells.Item(40, 10).Value = 400
$ws.Cells.Item(40, 11).Value = 7000
$ws.Cells.Item(40, 12).Value = 8000
$ws.Cells.Item(40, 13).Value = 7500
$ws.Cells.Item(40, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(40, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(40, 16).Value = 125
$ws.Cells.Item(40, 17).Value = 60
$ws.Cells.Item(40, 18).Value = "Hortaliza"
